# Applies targeted numeric-cell edits across ALC/ARM/BSM/GSM/LTW/WVR sheets
# (scheduled-runner profit recompute), per the authoritative diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 87
$ws.Cells.Item(87, 8).Value = 38583.332
$ws.Cells.Item(87, 9).Value = 30000
$ws.Cells.Item(87, 10).Value = 40300
$ws.Cells.Item(87, 11).Value = 30000
$ws.Cells.Item(87, 12).Value = 40300
$ws.Cells.Item(87, 13).Value = -28752
$ws.Cells.Item(87, 14).Value = -42796
# Row 90
$ws.Cells.Item(90, 8).Value = 38583.332
$ws.Cells.Item(90, 9).Value = 30000
$ws.Cells.Item(90, 10).Value = 40300
$ws.Cells.Item(90, 11).Value = 90000
$ws.Cells.Item(90, 12).Value = 120900
$ws.Cells.Item(90, 13).Value = -83760
$ws.Cells.Item(90, 14).Value = -133380

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Cells.Item(132, 8).Value = 2417.7778
$ws.Cells.Item(132, 9).Value = 1850.4286
$ws.Cells.Item(132, 11).Value = 5551.2858
$ws.Cells.Item(132, 13).Value = -3021.2858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 10160.826
$ws.Cells.Item(94, 9).Value = 1384.875
$ws.Cells.Item(94, 10).Value = 30220.143
$ws.Cells.Item(94, 11).Value = 1384.875
$ws.Cells.Item(94, 12).Value = 30220.143
$ws.Cells.Item(94, 13).Value = -933.875
$ws.Cells.Item(94, 14).Value = -31122.143

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Cells.Item(93, 8).Value = 8050.2
$ws.Cells.Item(93, 10).Value = 8050.2
$ws.Cells.Item(93, 12).Value = 8050.2
$ws.Cells.Item(93, 14).Value = -11794.2
# Row 122 (special case): H:L reset to 0, M/N cells removed entirely
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).ClearContents()
# Row 123
$ws.Cells.Item(123, 8).Value = 9721.182000000001
$ws.Cells.Item(123, 10).Value = 9721.182000000001
$ws.Cells.Item(123, 12).Value = 9721.182000000001
$ws.Cells.Item(123, 14).Value = -14621.182

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 123
$ws.Cells.Item(123, 8).Value = 32714.5
$ws.Cells.Item(123, 10).Value = 32714.5
$ws.Cells.Item(123, 12).Value = 32714.5
$ws.Cells.Item(123, 14).Value = -42514.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Cells.Item(33, 8).Value = 15110
$ws.Cells.Item(33, 9).Value = 15110
$ws.Cells.Item(33, 11).Value = 15110
$ws.Cells.Item(33, 13).Value = -14860
# Row 36
$ws.Cells.Item(36, 8).Value = 15110
$ws.Cells.Item(36, 9).Value = 15110
$ws.Cells.Item(36, 11).Value = 15110
$ws.Cells.Item(36, 13).Value = -14860
# Row 119
$ws.Cells.Item(119, 8).Value = 13449.667
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 13449.667
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 13449.667
$ws.Cells.Item(119, 14).Value = -23125.667
# Row 120
$ws.Cells.Item(120, 8).Value = 49500
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 49500
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 49500
$ws.Cells.Item(120, 14).Value = -59176
# Row 121
$ws.Cells.Item(121, 8).Value = 35420
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 35420
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 12).Value = 35420
$ws.Cells.Item(121, 14).Value = -38914
# Row 122
$ws.Cells.Item(122, 8).Value = 35715724
$ws.Cells.Item(122, 9).Value = 52632932
$ws.Cells.Item(122, 10).Value = 1618.3334
$ws.Cells.Item(122, 11).Value = 157898796
$ws.Cells.Item(122, 12).Value = 4855.0002
$ws.Cells.Item(122, 13).Value = -157896346
$ws.Cells.Item(122, 14).Value = -9755.0002
# Row 123
$ws.Cells.Item(123, 8).Value = 38000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 38000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 38000
$ws.Cells.Item(123, 14).Value = -47800
# Row 124
$ws.Cells.Item(124, 8).Value = 60899.8
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 60899.8
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 60899.8
$ws.Cells.Item(124, 14).Value = -70719.8
# Row 125
$ws.Cells.Item(125, 8).Value = 120000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 120000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 120000
$ws.Cells.Item(125, 14).Value = -129840
# Row 126
$ws.Cells.Item(126, 8).Value = 2239.1177
$ws.Cells.Item(126, 9).Value = 2474.3333
$ws.Cells.Item(126, 10).Value = 475
$ws.Cells.Item(126, 11).Value = 7422.999899999999
$ws.Cells.Item(126, 12).Value = 1425
$ws.Cells.Item(126, 13).Value = -4952.999899999999
$ws.Cells.Item(126, 14).Value = -6365
# Row 127
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
# Row 128
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
# Row 129
$ws.Cells.Item(129, 8).Value = 25000
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 25000
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 25000
$ws.Cells.Item(129, 14).Value = -35000
# Row 130
$ws.Cells.Item(130, 8).Value = 61964.5
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 61964.5
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 61964.5
$ws.Cells.Item(130, 14).Value = -72004.5
# Row 131
$ws.Cells.Item(131, 8).Value = 34250
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 34250
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 34250
$ws.Cells.Item(131, 14).Value = -44330
# Row 132
$ws.Cells.Item(132, 8).Value = 2580.8518
$ws.Cells.Item(132, 9).Value = 2762.2144
$ws.Cells.Item(132, 10).Value = 2385.5386
$ws.Cells.Item(132, 11).Value = 8286.643199999999
$ws.Cells.Item(132, 12).Value = 7156.6158
$ws.Cells.Item(132, 13).Value = -5756.643199999999
$ws.Cells.Item(132, 14).Value = -12216.6158
# Row 133
$ws.Cells.Item(133, 8).Value = 45443
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 45443
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 45443
$ws.Cells.Item(133, 14).Value = -55563
# Row 135
$ws.Cells.Item(135, 8).Value = 42900
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 42900
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 42900
$ws.Cells.Item(135, 14).Value = -53040
# Row 136
$ws.Cells.Item(136, 8).Value = 3787.5278
$ws.Cells.Item(136, 9).Value = 932.7931
$ws.Cells.Item(136, 10).Value = 15614.286
$ws.Cells.Item(136, 11).Value = 2798.3793
$ws.Cells.Item(136, 12).Value = 46842.858
$ws.Cells.Item(136, 13).Value = -248.3793000000001
$ws.Cells.Item(136, 14).Value = -51942.858
# Row 137
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
# Row 138
$ws.Cells.Item(138, 8).Value = 20000
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 20000
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 20000
$ws.Cells.Item(138, 14).Value = -30280
# Row 139
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
# Row 140
$ws.Cells.Item(140, 8).Value = 118400
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 118400
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 118400
$ws.Cells.Item(140, 14).Value = -128760
# Row 141
$ws.Cells.Item(141, 8).Value = 39871.668
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 39871.668
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 39871.668
$ws.Cells.Item(141, 14).Value = -50231.668
